# Daily refresh of the cryptocurrency price/volume table.
# Mirrors a scheduled GitHub Actions scrape: per-row Price (D) and
# Volume(1h) (E) cells get the latest snapshot; rows 45/46 (OKB vs
# VeChain) additionally swapped rank order, so Coin/Link/Price/Volume
# are rewritten for those two rows as a block.
#
# Price cells are forced to Text format before assignment: the source
# values are plain display strings (some are dotted-thousands like
# "58.810.72", some carry significant trailing zeros like "13.00"),
# and we must not let Excel auto-convert them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $text) {
    $rng = $ws.Range($rangeAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# --- Rows 45 & 46 swapped rank order: VeChain now outranks OKB ---
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D45" "0.0277"
$ws.Range("E45").Value = "  +7.78%  "

$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D46" "37.00"
$ws.Range("E46").Value = "  +0.50%  "

# --- Per-row Price / Volume(1h) refresh ---
Set-TextValue "D2" "58.810.72"
$ws.Range("E2").Value = "  +1.45%  "
Set-TextValue "D3" "3.178.83"
$ws.Range("E3").Value = "  +2.28%  "
$ws.Range("E4").Value = "  -0.04%  "
Set-TextValue "D5" "536.02"
Set-TextValue "D6" "143.16"
$ws.Range("E6").Value = "  +0.79%  "
$ws.Range("E7").Value = "  -0.01%  "
Set-TextValue "D8" "3.175.91"
$ws.Range("E8").Value = "  +2.22%  "
Set-TextValue "D9" "0.450"
$ws.Range("E9").Value = "  +2.70%  "
Set-TextValue "D10" "7.24"
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("E11").Value = "  +1.56%  "
Set-TextValue "D12" "0.401"
$ws.Range("E12").Value = "  +4.04%  "
Set-TextValue "D13" "3.725.83"
$ws.Range("E13").Value = "  +2.31%  "
$ws.Range("E14").Value = "  +2.86%  "
Set-TextValue "D15" "26.04"
$ws.Range("E15").Value = "  -1.02%  "
$ws.Range("E16").Value = "  +1.93%  "
Set-TextValue "D17" "58.836.90"
$ws.Range("E17").Value = "  +1.33%  "
Set-TextValue "D18" "3.181.73"
$ws.Range("E18").Value = "  +2.24%  "
Set-TextValue "D19" "6.21"
$ws.Range("E19").Value = "  +1.89%  "
Set-TextValue "D20" "13.00"
$ws.Range("E20").Value = "  +1.07%  "
Set-TextValue "D21" "8.12"
$ws.Range("E21").Value = "  +0.39%  "
Set-TextValue "D22" "359.07"
$ws.Range("E22").Value = "  +5.99%  "
$ws.Range("E23").Value = "  +0.07%  "
Set-TextValue "D24" "0.517"
$ws.Range("E24").Value = "  +2.08%  "
Set-TextValue "D25" "68.64"
$ws.Range("E25").Value = "  +3.70%  "
Set-TextValue "D26" "0.171"
$ws.Range("E26").Value = "  +1.92%  "
Set-TextValue "D27" "0.0₃0963"
$ws.Range("E27").Value = "  +5.54%  "
$ws.Range("E28").Value = "  -0.09%  "
Set-TextValue "D29" "7.58"
$ws.Range("E29").Value = "  +4.53%  "
Set-TextValue "D30" "6.59"
$ws.Range("E30").Value = "  -0.41%  "
$ws.Range("E31").Value = "  +0.04%  "
Set-TextValue "D32" "1.91"
$ws.Range("E32").Value = "  +2.61%  "
Set-TextValue "D33" "21.44"
$ws.Range("E33").Value = "  +2.04%  "
Set-TextValue "D34" "1.23"
$ws.Range("E34").Value = "  +1.68%  "
$ws.Range("E35").Value = "  +6.88%  "
Set-TextValue "D36" "158.02"
$ws.Range("E36").Value = "  +2.49%  "
Set-TextValue "D37" "6.27"
$ws.Range("E37").Value = "  +3.69%  "
$ws.Range("E38").Value = "  -1.70%  "
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("E40").Value = "  +13.82%  "
Set-TextValue "D41" "0.0679"
$ws.Range("E41").Value = "  +2.00%  "
$ws.Range("E42").Value = "  +4.28%  "
$ws.Range("E43").Value = "  +4.55%  "
Set-TextValue "D44" "3.217.35"
$ws.Range("E44").Value = "  +2.08%  "
Set-TextValue "D47" "2.348.62"
$ws.Range("E47").Value = "  +2.05%  "
$ws.Range("E48").Value = "  +0.01%  "
Set-TextValue "D49" "1.03"
$ws.Range("E49").Value = "  +6.59%  "
Set-TextValue "D50" "20.78"
$ws.Range("E50").Value = "  +0.32%  "
Set-TextValue "D51" "6.10"
$ws.Range("E51").Value = "  +1.65%  "
